{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\nlet target = null;\nfor (const p of paras.items) {\n  if (p.text.trim() === \"What is Science?\") { target = p; break; }\n}\nconst list = target.list;\nlist.setLevelBullet(0, Word.ListBullet.custom, \"\u2022\".charCodeAt(0), \"\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs(21)\n$lt = $p.Range.ListFormat.ListTemplate\n$lvl = $lt.ListLevels.Item(1)\n$lvl.Font.Name = \"Arial\"\n"}
